$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to the new (code-friendly) header names ---
$ws.Range("A1").Value = "StudentId"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "RegistrationYear"
$ws.Range("E1").Value = "YearOfStudies"
$ws.Range("G1").Value = "CourseId"
$ws.Range("H1").Value = "ProjectGrade"
$ws.Range("I1").Value = "ExamGrade"

# --- Semester column (F): replace numeric placeholder values with Spring/Fall ---
$ws.Range("F2").Value = "Spring"
$ws.Range("F3").Value = "Fall"
$ws.Range("F4").Value = "Fall"
$ws.Range("F5").Value = "Spring"
$ws.Range("F6").Value = "Fall"
$ws.Range("F7").Value = "Spring"
$ws.Range("F8").Value = "Spring"

# --- Strip stray leading spaces from the last three students' names ---
$ws.Range("B6").Value = "Name7"
$ws.Range("C6").Value = "LName1"
$ws.Range("B7").Value = "Name8"
$ws.Range("C7").Value = "LName2"
$ws.Range("B8").Value = "Name9"
$ws.Range("C8").Value = "LName3"

# --- Rename the sheet-scoped defined name (table got re-imported as _1) ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!myy301_students") {
        $n.Name = "myy301_students_1"
    }
}

# --- Selection moved from D11 to D12 ---
[void]$ws.Range("D12").Select()
